# Populate row 2 (ICBM MGS Trainer Repair - replaces the old TMD Trophies row content)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = @'
ICBM MGS Trainer Repair
'@
$ws.Range("B2").Value = @'
https://sam.gov/opp/025035bc55484bc4a2369385aa613f01/view
'@
$ws.Range("C2").Value = @'
ICBM MGS Trainer Repair
Notice ID: FA820622R0006

Solicitation amended to:
1. Request proposal for FIVE (5) YEAR IDIQ TYPE CONTRACT FOR REPAIR OF M

...

Department/Ind.Agency
DEPT OF DEFENSE
Subtier
DEPT OF THE AIR FORCE
Office
FA8206 AFSC PZAAA2
'@
$ws.Range("D2").Value = @'
Contract Opportunities
Current Date Offers Due
July 14, 2022, 11:00 PM GMT+1
Notice Type
Updated Solicitation
Updated Date
Jun 6, 2022 (1)
Published Date
Jun 6, 2022
'@
$ws.Range("F2").Value = ' Original Set Aside: '
$ws.Range("G2").Value = @'
Contract Opportunity General Information  Classification  Description  Attachments/Links  Contact Information  History  ICBM MGS Trainer Repair Active Contract Opportunity  FA820622R0006  Related Notice  Contract Line Item Number  Department/Ind. Agency  General Information Contract Opportunity Type:  Solicitation  (Updated)All Dates/Times are:  (UTC-06:00) MOUNTAIN STANDARD TIME, DENVER, USA Updated Published Date:  Jun 06, 2022 12:58 pm MDT Original Published Date: Inactive Policy: 15 days after date offers due Updated Inactive Date:  Jul 29, 2022 Original Inactive Date: Authority: Fair Opportunity / Limited Sources Justification Authority: Initiative: None Classification  Original Set Aside: Product Service Code: NAICS Code: Place of Performance:          USA DescriptionSolicitation amended to: 
1. Request proposal for FIVE (5) YEAR IDIQ TYPE CONTRACT FOR REPAIR OF MSG TRAINERS
This is a single award Indefinite Delivery/Indefinite Quantity (IDIQ), Firm Fixed Priced solicitation for a five year ordering period.
Quantities ordered shall be determined at the issuance of each Task Order/Delivery Order.
       The minimum order amount is repair of 3 MGS trainers.
       The estimated maximum amount is repair of 5 MGS trainers.
2. Add Over and Above CLINs added to be negotiated at each task order.
3. Extend receipt of proposal request to 14 July 2022.
4. Attach updated PWS to reflect the min and max quantity, and address over and above process.
Additional information request:
With proposal response, please identify what aspects of repair would be included in the basic service, priced under CLIN X002, and what types of repair would be over and above that proposed price, which may apply to CLIN X009.
Original solicitation information: 
This is a solicitation notice of intent to award a sole source contract in support of the ICBM Minuteman II (MMIII) Trainer Missile Guidance Sets (MGS Trainer). Air Force Materiel Command, Air Force Sustainment Center, Hill AFB, UT, intends to award a Sole Source FFP contract as implemented by the Federal Acquisition Regulation Subpart 15.2, as supplemented with additional information in this notice. The contract will be awarded to The Boeing Company. IAW FAR 6.302-1, Market research was previously conducted and it has been determined that Boeing is the sole source that can meet the Government's requirements.
This is not a request for competitive proposals and no solicitation is available. All responsible sources may submit a capability statement which shall be considered by the agency. Interested sources must respond in writing with clear and convincing evidence to support their ability to provide the required supplies/services within fifteen (15) days of this publication. A determination by the Government not to compete this action based on responses received is solely within the discretion of the Government.
Contact Information  Primary Point of Contact Terra Clarke     terra.clarke@us.af.mil    Phone Number 8015863954  Secondary Point of Contact Paige LaPoint     paige.lapoint@us.af.mil    Phone Number 8015862947 History
'@
$ws.Range("H2").Value = 'IDIQ'
# --- Row 3 ---
$ws.Range("A3").Value = @'
Removal and Replacement of AHU-1 Chilled Water Coils - NLAE
'@
$ws.Range("B3").Value = @'
https://sam.gov/opp/4c74d299b5054a14af636f93b53b0f99/view
'@
$ws.Range("C3").Value = @'
Removal and Replacement of AHU-1 Chilled Water Coils - NLAE
Notice ID: 12505B22R0016

The United States Department of Agriculture (USDA), Agricultural Research Service (ARS) has a requirement for Replace AHU-1 chiller water coils.

...

Department/Ind.Agency
AGRICULTURE, DEPARTMENT OF
Subtier
AGRICULTURAL RESEARCH SERVICE
Office
USDA ARS MWA AAO ACQ/PER PROP
'@
$ws.Range("D3").Value = @'
Contract Opportunities
Current Date Offers Due
July 13, 2022, 04:00 PM GMT+1
Notice Type
Original Solicitation
Updated Date
Jun 6, 2022
Published Date
Jun 6, 2022
'@
$ws.Range("F3").Value = ' Original Set Aside: '
$ws.Range("G3").Value = @'
Contract Opportunity General Information  Classification  Description  Attachments/Links  Contact Information  History  Removal and Replacement of AHU-1 Chilled Water Coils - NLAE Active Contract Opportunity  Notice ID  12505B22R0016  Related Notice  Department/Ind. Agency  General Information Contract Opportunity Type: Solicitation (Original) All Dates/Times are:  (UTC-05:00) CENTRAL STANDARD TIME, CHICAGO, USA Original Published Date: Jun 06, 2022 01:56 pm CDTOriginal Date Offers Due: Jul 13, 2022 10:00 am CDTInactive Policy: 15 days after date offers due Original Inactive Date: Jul 28, 2022Initiative: None Classification  Original Set Aside: Product Service Code: NAICS Code: Place of Performance:     ,   50011  USA DescriptionThe United States Department of Agriculture (USDA), Agricultural Research Service (ARS) has a requirement for Replace AHU-1 chiller water coils.
Solicitation 12505B22R0016 is issued as a Request for Proposals (RFP).  The solicitation document, with the incorporated clauses and provisions, is issued in accordance with the Federal Acquisition Regulation (FAR) Parts 13 and 36 and will be incorporated in the awarded firm-fixed price contract.  Interested vendors should reference the solicitation documents for all terms and conditions.  All interested parties are responsible for monitoring this website to ensure they have the most current information (i.e. Amendments) for the solicitation.  No paper solicitation will be available.
This requirement is solicited as unrestricted for full and open competition.
The associated NAICS Code is 238220, Plumbing, Heating, and Air-Conditioning Contractors.  The small business size standard is $16.50 million.
Magnitude of Construction is estimated to be between $25,000 and $100,000.
Period of performance is 59 days after receipt of the Notice to Proceed to include inspection and punch list.
Site Visit:  Organized site visits are scheduled for June 14, 2022 at 10:00 AM CT and June 15, 2022 at 1:00 PM CT.  Face coverings must be worn at all times.  The COVID-19 Vaccine Attestation form must be presented upon entry to the site visit.  The site visit will be held at the project location.  The same information will be provided at all site visits, so offerors only need to attend one.  In order to attend the site visit, you MUST pre-register by 2:00 PM CST on June 13, 2022.  Pre-register by emailing melissa.grice@usda.gov with your company name, the name(s) of the representatives who will attend, and which site visit you are registering for.  An email response will be sent confirming your registration. 
The site visit will be held for the purpose of providing contractors with the opportunity to familiarize themselves with the site which may be helpful in the preparation of offers.  Attendance at this site visit is not mandatory for offer submission; however, failure to visit the site will not relieve or mitigate the successful contractor's responsibility and obligation to fully comply with the terms, conditions, and specifications contained and/or referenced in this document.
All questions regarding this solicitation must be submitted in writing to the Contract Specialist, Melissa Grice, via email to melissa.grice@usda.gov.  Questions must be submitted no later than close of business on July 1, 2022.  Answers to all questions received by that time will be posted as an amendment to the solicitation.  No questions will be answered after this date unless determined to be in the best interest of the Government as deemed by the Contracting Officer.  Telephone requests for information will not be accepted or returned.
Interested offerors must be registered in the System for Award Management (SAM).  To register go to www.sam.gov.  Instructions for registering are on the web page (there is no fee for registration).
The solicitation and associated information will be available only from the Contract Opportunities page at SAM.gov.
This solicitation is NOT an invitation for bids and there will be NO formal public bid opening.  All inquiries must be in writing via email to the persons specified in this solicitation.  All answers will be provided in writing via posting to the web. 
DISCLAIMER: The official contract documents are located on the government webpage and the Government is not liable for information furnished by any other source.  Amendments, if/when issued will be posted to Contract Opportunities for electronic downloading.  This will normally be the only method of distributing amendments prior to closing; therefore, it is the offeror’s responsibility to check the website periodically for any amendments to this solicitation.  Websites are occasionally inaccessible, due to various reasons.  The Government is not responsible for any loss of internet connectivity or for an offeror’s inability to access the documents posted on the referenced web pages.  The Government will not issue paper copies.
Biobased Products:  This procurement requires the use of Biobased products to the extent that such products are reasonably available, meet agency or relevant industry performance standards, and are reasonably priced.  The products should first be acquired from the USDA designated product categories.  All supplies and materials shall be of a type and quantity that conform to applicable Federal specifications and standards, and to the extent feasible and reasonable, include the exclusive use of biobased and recycled products.  Please visit www.biopreferred.gov for more information on the USDA Biobased Program and to reference the catalog of mandatory biobased products.
Attachments/LinksContact Information  Primary Point of Contact Missy Grice     melissa.grice@usda.gov    Phone Number 9706314559  Secondary Point of Contact History
'@
$ws.Range("H3").Value = @'
BPA
'@
# --- Row 4 ---
$ws.Range("A4").Value = @'
300WM Ammunition
'@
$ws.Range("B4").Value = @'
https://sam.gov/opp/1965368aa94347229ec5d5111db2da7b/view
'@
$ws.Range("C4").Value = @'
300WM Ammunition
Notice ID: 70US0922Q70092047

**Questions & Answers - attached as of 6 June 2022**
The purpose of this solicitation is to establish and Ind

...

Department/Ind.Agency
HOMELAND SECURITY, DEPARTMENT OF
Subtier
U.S. SECRET SERVICE
Office
U. S. SECRET SERVICE
'@
$ws.Range("D4").Value = @'
Contract Opportunities
Current Date Offers Due
July 15, 2022, 11:20 PM GMT+1
Notice Type
Updated Solicitation
Updated Date
Jun 6, 2022 (1)
Published Date
Jun 6, 2022
'@
$ws.Range("F4").Value = ' Original Set Aside: '
$ws.Range("G4").Value = @'
Contract Opportunity General Information  Classification  Description  Attachments/Links  Contact Information  History  300WM Ammunition Active Contract Opportunity  70US0922Q70092047  Related Notice  Contract Line Item Number  Department/Ind. Agency  General Information Contract Opportunity Type:  Solicitation  (Updated)All Dates/Times are:  (UTC-04:00) EASTERN STANDARD TIME, NEW YORK, USA Updated Published Date:  Jun 06, 2022 12:18 pm EDT Original Published Date: Inactive Policy: 15 days after date offers due Updated Inactive Date:  Jul 30, 2022 Original Inactive Date: Authority: Fair Opportunity / Limited Sources Justification Authority: Initiative: None Classification  Original Set Aside: Product Service Code: NAICS Code: Place of Performance:     ,   20708  USA Description**Questions & Answers - attached as of 6 June 2022**
The purpose of this solicitation is to establish and Indefinite Delivery/Indefinite Quantity (IDIQ) for 300WM  ammunition contract for the United States Secret Service (USSS), Roqley Training Center (RTC) in support of firearms training and duty.
Contact Information  Primary Point of Contact Crystal Pressley     crystal.pressley@usss.dhs.gov    Phone Number 2024065209  Secondary Point of Contact Theresa Williams     theresa.williams@usss.dhs.gov    Phone Number 2024066213 History
'@
$ws.Range("H4").Value = 'IDIQ'

# Re-fit row heights so they match the default (no custom height), since the
# multi-line cell values above would otherwise trigger Excel auto row-height.
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()

